# Quarterly database update for شپدیس (Shapdis) income statement.
# The spreadsheet keeps a rolling window of the 10 most-recent quarters
# in columns D:M. This edit rolls the window forward by one quarter:
# every existing column's data moves one column to the left (D<-E,
# E<-F, ... L<-M) and the brand-new quarter's figures are written into
# the now-vacant column M. The same left-shift applies to the period
# labels (row 8) and the publish-date labels (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M")

# ---------------------------------------------------------------------
# Row 8: financial period headers (shift left, new quarter appended)
# ---------------------------------------------------------------------
$periods = @(
  "فصل اول منتهی به 1399/09",
  "فصل دوم منتهی به 1399/12",
  "فصل سوم منتهی به 1400/03",
  "فصل چهارم منتهی به 1400/06",
  "فصل اول منتهی به 1400/09",
  "فصل دوم منتهی به 1400/12",
  "فصل سوم منتهی به 1401/03",
  "فصل چهارم منتهی به 1401/06",
  "فصل اول منتهی به 1401/09",
  "فصل دوم منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
}

# ---------------------------------------------------------------------
# Row 9: publish-date headers (shift left, new quarter appended)
# ---------------------------------------------------------------------
$dates = @(
  "1400-10-29 (2)",
  "1401-03-21 (4)",
  "1401-05-19 (3)",
  "1401-10-05 (9)",
  "1401-10-28 (2)",
  "1402-01-30 (3)",
  "1401-05-19 (2)",
  "1402-01-30 (4)",
  "1401-10-28",
  "1402-01-30"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $dates[$i]
}

# ---------------------------------------------------------------------
# Data rows 11-25, 27: shift left, new quarter's values appended in M
# ---------------------------------------------------------------------
$rowData = @{
  11 = @(46810537, 41995653, 47934192, 60658299, 94387382, 66986462, 128059977, 142149002, 110842812, 70344021)
  12 = @(-14466354, -15660374, -14087457, -26496392, -37794904, -31128186, -41989284, -48770181, -34318064, -30885267)
  13 = @(32344183, 26335279, 33846735, 34161907, 56592478, 35858276, 86070693, 93378821, 76524748, 39458754)
  14 = @(-3831223, -3825327, -3138203, 1202914, -2265379, -1341498, -1672323, -33567943, -3191115, -43666810)
  16 = @(273315, 3665207, -2373168, 807037, -169830, 1637528, -440850, 10227527, 1836195, 21479473)
  17 = @(28786275, 26175159, 28335364, 36171858, 54157269, 36154306, 83957520, 70038405, 75169828, 17271417)
  18 = @(-77409, -172357, -66393, -157919, -258982, -179286, -149964, -477395, -484883, -469945)
  19 = @(115810, 1031846, 1201308, 2176262, 2410362, 3870256, 2820935, 3669171, 3330336, 12255816)
  20 = @(28824676, 27034648, 29470279, 38190201, 56308649, 39845276, 86628491, 73230181, 78015281, 29057288)
  21 = @(0, -26316, 26316, 4597514, 0, -16130, 16130, -184701, 0, -41198)
  22 = @(28824676, 27008332, 29496595, 42787715, 56308649, 39829146, 86644621, 73045480, 78015281, 29016090)
  24 = @(28824676, 27008332, 29496595, 42787715, 56308649, 39829146, 86644621, 73045480, 78015281, 29016090)
  25 = @(4804, 4501, 4916, 7131, 9385, 6638, 14441, 12174, 13003, 4836)
  27 = @(4804, 4501, 4916, 7131, 9385, 6638, 14441, 12174, 13003, 4836)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# Row 26 (سرمایه / capital): values are unchanged except a data fix in
# column J where a stray reported figure is corrected back to 6,000,000.
$ws.Range("J26").Value = 6000000

# ---------------------------------------------------------------------
# Column widths: the rolling shift also moves the "wide" (31) columns
# that mark the boundary of each year from D/H/L to G/K (i.e. the
# widths shift left by one column, same as the data).
# ---------------------------------------------------------------------
$narrow = 28.166666666666668   # -> serializes as width 29
$wide   = 30.166666666666668   # -> serializes as width 31

$ws.Columns.Item(4).ColumnWidth  = $narrow   # D
$ws.Columns.Item(5).ColumnWidth  = $narrow   # E
$ws.Columns.Item(6).ColumnWidth  = $narrow   # F
$ws.Columns.Item(7).ColumnWidth  = $wide     # G
$ws.Columns.Item(8).ColumnWidth  = $narrow   # H
$ws.Columns.Item(9).ColumnWidth  = $narrow   # I
$ws.Columns.Item(10).ColumnWidth = $narrow   # J
$ws.Columns.Item(11).ColumnWidth = $wide     # K
$ws.Columns.Item(12).ColumnWidth = $narrow   # L
$ws.Columns.Item(13).ColumnWidth = $narrow   # M

# ---------------------------------------------------------------------
# Row heights: a handful of rows carry an explicit custom height; keep
# them in sync with the refreshed layout.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6
